$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.947.34"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.210.99"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'240.73"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "'72.40"
$ws.Range("E7").Value = "  -4.55%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  -3.79%  "
$ws.Range("D10").Value = "'42.00"
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "2.547.16"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").Value = "'14.17"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "'0.831"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "2.209.72"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "41.824.97"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'0.0000105"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "'72.46"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'6.13"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'10.56"
$ws.Range("E22").Value = "  +15.67%  "
$ws.Range("D23").Value = "'228.69"
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  -7.95%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'11.43"
$ws.Range("D27").Value = "'3.63"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D30").Value = "'167.24"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'20.40"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("E32").Value = "  +5.33%  "
$ws.Range("D33").Value = "'0.0790"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").Value = "'29.93"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "  -12.79%  "
$ws.Range("D37").Value = "'4.19"
$ws.Range("E37").Value = "  -7.09%  "
$ws.Range("D38").Value = "'0.0298"
$ws.Range("E38").Value = "  -6.19%  "
$ws.Range("D39").Value = "'13.76"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "'64.54"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").Value = "'2.10"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").Value = "'5.62"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("D43").Value = "'0.195"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("D44").Value = "'8.65"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "'103.71"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "'2.34"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "2.420.38"
$ws.Range("E51").Value = "  -2.01%  "
